# Update "想去人数" (want-to-go count) figures in column F
# for the "展览" and "全部类型" sheets, which carry duplicate data.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 1743
    4  = 799
    5  = 3
    8  = 12095
    9  = 43
    13 = 1118
    15 = 13521
    16 = 13570
    18 = 157
    21 = 998
    24 = 2028
    25 = 188
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
